# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp.
# - Refresh case counts for several countries (India, Israel, Armenia,
#   Afganistan, Hungria, Georgia).
# - Insert a fresh "El Salvador" row (with updated stats) right after
#   "Costa Rica", and remove its old row (which used to sit between
#   "Australia" and "Etiopia").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1, column A) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 09:15"

# --- Refresh stats for countries whose figures simply changed ---

# India (row 6)
$ws.Range("B6").Value = 1535335
$ws.Range("C6").Value = 3200
$ws.Range("D6").Value = 989878
$ws.Range("E6").Value = 511205
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 28
$ws.Range("H6").Value = 34252

# Israel (row 39)
$ws.Range("B39").Value = 66555
$ws.Range("C39").Value = 262
$ws.Range("D39").Value = 32692
$ws.Range("E39").Value = 33377

# Armenia (row 53)
$ws.Range("B53").Value = 37937
$ws.Range("C53").Value = 308
$ws.Range("D53").Value = 27824
$ws.Range("E53").Value = 9390
$ws.Range("G53").Value = 4
$ws.Range("H53").Value = 723

# Afganistan (row 54)
$ws.Range("B54").Value = 36471
$ws.Range("C54").Value = 103
$ws.Range("D54").Value = 25389
$ws.Range("E54").Value = 9811
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 1271

# Hungria (row 103)
$ws.Range("B103").Value = 4465
$ws.Range("C103").Value = 9
$ws.Range("D103").Value = 3339
$ws.Range("E103").Value = 530

# Georgia (row 143)
$ws.Range("B143").Value = 1155
$ws.Range("C143").Value = 10
$ws.Range("D143").Value = 929
$ws.Range("E143").Value = 210

# --- Move / refresh "El Salvador" ---
# Insert a new row right after "Costa Rica" (row 71) for the updated
# "El Salvador" entry, pushing Chequia / Costa de Marfil / Australia down.
$ws.Rows(72).Insert()
$ws.Range("A72").Value = "El Salvador"
$ws.Range("B72").Value = 15841
$ws.Range("C72").Value = 395
$ws.Range("D72").Value = 8071
$ws.Range("E72").Value = 7340
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 13
$ws.Range("H72").Value = 430

# Remove the old "El Salvador" row, now shifted down to row 76 (it used
# to sit between Australia and Etiopia).
$ws.Rows(76).Delete()
